$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal text, preventing Excel's
# automatic conversion of date-looking strings ("11/20/2022", etc.)
# into date serial numbers.
function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Update the release schedule dates for the Forecasting and Machine
# Learning sections.
Set-TextValue "C12" "11/20/2022"
Set-TextValue "C13" "11/20/2022"
Set-TextValue "C14" "11/27/2022"
Set-TextValue "C16" "11/13/2022"
Set-TextValue "C25" "11/06/2022"

# Reflect the saved selection state.
$ws.Range("C15").Select()
